$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column N holds literal date-like text (e.g. "2022-11-15"), not real dates.
# Force text format first so Excel does not auto-convert the string into a
# date serial number, then restore the default "Normal" style so the saved
# cell has no lingering number-format override (matches the original cells).
$ws.Range("N2:N18").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "HALIFAX, NS"
$ws.Range("E2").Value = "MAERSK IDAHO"
$ws.Range("F2").Value = "MKIDO"
$ws.Range("G2").Value = "0EDCME1MA"
$ws.Range("H2").Value = "MSK"
$ws.Range("I2").Value = 9193264
$ws.Range("J2").Value = 44880.83333333334
$ws.Range("K2").Value = 44881.375
$ws.Range("L2").Value = 44881.4166666551
$ws.Range("M2").Value = 44881.9166666551
$ws.Range("N2").Value = "2022-11-15"
$ws.Range("O2").Value = "CAHAL"

# Row 3
$ws.Range("D3").Value = "MIAMI"
$ws.Range("E3").Value = "GARFIELD"
$ws.Range("F3").Value = "GAFLD"
$ws.Range("G3").Value = "0GY9EE1MA"
$ws.Range("H3").Value = "CMA - CGM"
$ws.Range("I3").Value = 9311842
$ws.Range("J3").Value = 44880.875
$ws.Range("K3").Value = 44882.0416666551
$ws.Range("L3").Value = 44882.125
$ws.Range("M3").Value = 44882.6666666551
$ws.Range("N3").Value = "2022-11-15"
$ws.Range("O3").Value = "USMIA"

# Row 4
$ws.Range("D4").Value = "NEW ORLEANS"
$ws.Range("E4").Value = "CMA CGM VERACRUZ"
$ws.Range("F4").Value = "CGVCZ"
$ws.Range("G4").Value = "0GBDMS1MA"
$ws.Range("I4").Value = 9418377
$ws.Range("J4").Value = 44880.9166666551
$ws.Range("K4").Value = 44883.5416666551
$ws.Range("L4").Value = 44884.0416666551
$ws.Range("M4").Value = 44884.7916666551
$ws.Range("N4").Value = "2022-11-15"
$ws.Range("O4").Value = "USMSY"

# Row 5
$ws.Range("D5").Value = "PORT HUENEME, CA"
$ws.Range("E5").Value = "DEL MONTE HARVESTER"
$ws.Range("F5").Value = "DELHV"
$ws.Range("G5").Value = "0WS1LR1MA"
$ws.Range("H5").Value = "DEL MONTE FRESH COMPANY"
$ws.Range("I5").Value = 9869667
$ws.Range("J5").Value = 44880.83333333334
$ws.Range("K5").Value = 44883.34861111111
$ws.Range("L5").Value = 44883.3902777662
$ws.Range("M5").Value = 44884.93194444444
$ws.Range("N5").Value = "2022-11-15"
$ws.Range("O5").Value = "USNTD"

# Row 6
$ws.Range("D6").Value = "NEW YORK"
$ws.Range("E6").Value = "NORTHERN PRIORITY"
$ws.Range("F6").Value = "NTPRY"
$ws.Range("G6").Value = "0AMCFS1MA"
$ws.Range("H6").Value = "MSK"
$ws.Range("I6").Value = 9450313
$ws.Range("J6").Value = 44880.875
$ws.Range("K6").Value = 44882.33333333334
$ws.Range("L6").Value = 44882.5
$ws.Range("M6").Value = 44883.1666666551
$ws.Range("N6").Value = "2022-11-15"
$ws.Range("O6").Value = "USNYC"

# Row 7
$ws.Range("D7").Value = "NEW YORK"
$ws.Range("E7").Value = "POLYNESIA"
$ws.Range("F7").Value = "PLYSA"
$ws.Range("G7").Value = "0UAD1S1MA"
$ws.Range("H7").Value = "CMA - CGM"
$ws.Range("I7").Value = 9477347
$ws.Range("J7").Value = 44880.875
$ws.Range("K7").Value = 44884.375
$ws.Range("L7").Value = 44884.5416666551
$ws.Range("M7").Value = 44885.375
$ws.Range("N7").Value = "2022-11-15"
$ws.Range("O7").Value = "USNYC"

# Row 8
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = "NEW YORK"
$ws.Range("E8").Value = "COSCO FORTUNE"
$ws.Range("F8").Value = "COFOR"
$ws.Range("G8").Value = "0MBCEW1MA"
$ws.Range("I8").Value = 9472127
$ws.Range("J8").Value = 44880.875
$ws.Range("K8").Value = 44876.3125
$ws.Range("L8").Value = 44880.45833333334
$ws.Range("M8").Value = 44883.75555555556
$ws.Range("N8").Value = "2022-11-15"
$ws.Range("O8").Value = "USNYC"

# Row 9
$ws.Range("D9").Value = "NEW YORK"
$ws.Range("E9").Value = "EVER FASHION"
$ws.Range("F9").Value = "EVFSH"
$ws.Range("G9").Value = "0VCD2W1MA"
$ws.Range("H9").Value = "EVERGREEN MARINE CORPORATION LTD"
$ws.Range("I9").Value = 9850836
$ws.Range("J9").Value = 44880.875
$ws.Range("K9").Value = 44878.33333333334
$ws.Range("L9").Value = 44881.45833333334
$ws.Range("M9").Value = 44883.2916666551
$ws.Range("N9").Value = "2022-11-15"
$ws.Range("O9").Value = "USNYC"

# Row 10
$ws.Range("D10").Value = "OAKLAND, CA"
$ws.Range("E10").Value = "EVER LINKING"
$ws.Range("F10").Value = "EVLKG"
$ws.Range("G10").Value = "0TBDGW1MA"
$ws.Range("H10").Value = "EVERGREEN MARINE CORPORATION LTD"
$ws.Range("I10").Value = 9629043
$ws.Range("J10").Value = 44880
$ws.Range("K10").Value = 44880.9166666551
$ws.Range("L10").Value = 44881.0416666551
$ws.Range("M10").Value = 44883.08333333334
$ws.Range("N10").Value = "2022-11-15"
$ws.Range("O10").Value = "USOAK"

# Row 11
$ws.Range("D11").Value = "OAKLAND, CA"
$ws.Range("E11").Value = "CMA CGM NEW JERSEY"
$ws.Range("F11").Value = "CGJEY"
$ws.Range("G11").Value = "0GVCDW1MA"
$ws.Range("H11").Value = "CMA - CGM"
$ws.Range("I11").Value = 9351141
$ws.Range("J11").Value = 44880
$ws.Range("K11").Value = 44881.9166666551
$ws.Range("L11").Value = 44882
$ws.Range("M11").Value = 44883.5416666551
$ws.Range("N11").Value = "2022-11-15"
$ws.Range("O11").Value = "USOAK"

# Row 12
$ws.Range("D12").Value = "NORFOLK"
$ws.Range("E12").Value = "MAERSK CHICAGO"
$ws.Range("F12").Value = "MKCHG"
$ws.Range("G12").Value = "1JU1ZE1MA"
$ws.Range("H12").Value = "MSK"
$ws.Range("I12").Value = 9332975
$ws.Range("J12").Value = 44880.875
$ws.Range("K12").Value = 44881.375
$ws.Range("L12").Value = 44881.5416666551
$ws.Range("M12").Value = 44881.95833333334
$ws.Range("N12").Value = "2022-11-15"
$ws.Range("O12").Value = "USORF"

# Row 13
$ws.Range("D13").Value = "SAVANNAH"
$ws.Range("E13").Value = "CMA CGM LA TRAVIATA"
$ws.Range("F13").Value = "CMTRA"
$ws.Range("G13").Value = "0MRBUE1MA"
$ws.Range("H13").Value = "CMA - CGM"
$ws.Range("I13").Value = 9299795
$ws.Range("J13").Value = 44880.875
$ws.Range("K13").Value = 44881.69097222222
$ws.Range("L13").Value = 44882.45833333334
$ws.Range("M13").Value = 44883.95833333334
$ws.Range("N13").Value = "2022-11-15"
$ws.Range("O13").Value = "USSAV"

# Row 14
$ws.Range("D14").Value = "SAVANNAH"
$ws.Range("E14").Value = "CMA CGM APOLLON"
$ws.Range("F14").Value = "CGAPO"
$ws.Range("G14").Value = "0MBC8W1MA"
$ws.Range("I14").Value = 9882516
$ws.Range("J14").Value = 44880.83333333334
$ws.Range("K14").Value = 44866.95833333334
$ws.Range("L14").Value = 44883.2916666551
$ws.Range("M14").Value = 44885.95833333334
$ws.Range("N14").Value = "2022-11-15"
$ws.Range("O14").Value = "USSAV"

# Row 15
$ws.Range("D15").Value = "SAVANNAH"
$ws.Range("E15").Value = "SEATRADE BLUE"
$ws.Range("F15").Value = "STBLE"
$ws.Range("G15").Value = "0RPBEN1MA"
$ws.Range("I15").Value = 9756107
$ws.Range("J15").Value = 44880.875
$ws.Range("K15").Value = 44881.375
$ws.Range("L15").Value = 44881.5416666551
$ws.Range("M15").Value = 44882.1666666551
$ws.Range("N15").Value = "2022-11-15"
$ws.Range("O15").Value = "USSAV"

# Row 16
$ws.Range("D16").Value = "TACOMA, WA"
$ws.Range("E16").Value = "EVER SUMMIT"
$ws.Range("F16").Value = "EVSUM"
$ws.Range("G16").Value = "0NWD2W1MA"
$ws.Range("H16").Value = "EVERGREEN MARINE CORPORATION LTD"
$ws.Range("I16").Value = 9300453
$ws.Range("J16").Value = 44880
$ws.Range("K16").Value = 44881.2916666551
$ws.Range("L16").Value = 44881.5416666551
$ws.Range("M16").Value = 44883.5416666551
$ws.Range("N16").Value = "2022-11-15"
$ws.Range("O16").Value = "USTIW"

# Row 17
$ws.Range("D17").Value = "PORT EVERGLADES"
$ws.Range("E17").Value = "POLAR CHILE"
$ws.Range("F17").Value = "POCHL"
$ws.Range("G17").Value = "0AMCCN1MA"
$ws.Range("H17").Value = "MSK"
$ws.Range("I17").Value = 9797187
$ws.Range("J17").Value = 44880.875
$ws.Range("K17").Value = 44886.6666666551
$ws.Range("L17").Value = 44886.75
$ws.Range("M17").Value = 44887.125
$ws.Range("N17").Value = "2022-11-15"
$ws.Range("O17").Value = "USPEF"

# Row 18
$ws.Range("D18").Value = "PORT EVERGLADES"
$ws.Range("E18").Value = "NYK RUMINA"
$ws.Range("F18").Value = "NYRNA"
$ws.Range("G18").Value = "0CLCPW1MA"
$ws.Range("H18").Value = "OCEAN NETWORK EXPRESS PTE. LTD."
$ws.Range("I18").Value = 9416991
$ws.Range("J18").Value = 44880.8541666551
$ws.Range("K18").Value = 44882.45833333334
$ws.Range("L18").Value = 44882.5
$ws.Range("M18").Value = 44883.1666666551
$ws.Range("N18").Value = "2022-11-15"
$ws.Range("O18").Value = "USPEF"

$ws.Range("N2:N18").Style = "Normal"

# Remove now-obsolete rows 19-23 (data consolidated down to 18 rows)
$ws.Rows("19:23").Delete()

